$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The cell that held the bare filename "jose_miguel_gutierrez.jpg" now holds
# the full raw-GitHub URL for the photo, and becomes a clickable hyperlink.
$url = "https://raw.githubusercontent.com/Pavanona/Diputados/refs/heads/main/jose_miguel_gutierrez.jpg"

$cell = $ws.Range("X2")
$cell.Value = $url
$cell.Hyperlinks.Add($cell, $url)

# The selection left active after the edit moved down one row.
$ws.Range("X3").Select()
